$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("M3").Value = 1.13
$ws.Range("N3").Value = 6
$ws.Range("Q3").Value = 2.88
$ws.Range("R3").Value = 1.4

# Row 6
$ws.Range("G6").Value = 1.85
$ws.Range("H6").Value = 3.3
$ws.Range("I6").Value = 4.05
$ws.Range("J6").Value = 2.4
$ws.Range("K6").Value = 2.1
$ws.Range("L6").Value = 4.45
$ws.Range("O6").Value = 1.34
$ws.Range("P6").Value = 2.75
$ws.Range("T6").Value = 2.55
$ws.Range("U6").Value = 1.85
$ws.Range("X6").Value = 8.25
$ws.Range("Z6").Value = 15.5
$ws.Range("AA6").Value = 15.5
$ws.Range("AC6").Value = 8.5
$ws.Range("AG6").Value = 10.25
$ws.Range("AH6").Value = 22
$ws.Range("AI6").Value = 14
$ws.Range("AJ6").Value = 65
$ws.Range("AK6").Value = 45
$ws.Range("AL6").Value = 50
$ws.Range("AN6").Value = 3.6
$ws.Range("AO6").Value = 9
$ws.Range("AP6").Value = 18
$ws.Range("AQ6").Value = 32
$ws.Range("AR6").Value = 65
$ws.Range("AT6").Value = 2.5
$ws.Range("AW6").Value = 5.8
$ws.Range("AX6").Value = 23
$ws.Range("AY6").Value = 30
$ws.Range("AZ6").Value = 150
$ws.Range("BA6").Value = 175
$ws.Range("BB6").Value = 400

# Row 7
$ws.Range("G7").Value = 2.05
$ws.Range("I7").Value = 3.4
$ws.Range("J7").Value = 2.63
$ws.Range("L7").Value = 3.75
$ws.Range("S7").Value = 1.3
$ws.Range("T7").Value = 3.4
$ws.Range("W7").Value = 10
$ws.Range("Y7").Value = 9
$ws.Range("Z7").Value = 19
$ws.Range("AD7").Value = 7.5
$ws.Range("AE7").Value = 12
$ws.Range("AH7").Value = 19
$ws.Range("AI7").Value = 12
$ws.Range("AJ7").Value = 34
$ws.Range("AK7").Value = 23
$ws.Range("AM7").Value = 126
$ws.Range("AO7").Value = 11
$ws.Range("AT7").Value = 3.4
$ws.Range("AX7").Value = 17
$ws.Range("BB7").Value = 126

# Row 10
$ws.Range("G10").Value = 3.6
$ws.Range("I10").Value = 1.95
$ws.Range("J10").Value = 4.5
$ws.Range("L10").Value = 2.75
$ws.Range("N10").Value = 8.5
$ws.Range("W10").Value = 8.5
$ws.Range("X10").Value = 17
$ws.Range("AH10").Value = 8.5
$ws.Range("AJ10").Value = 17
$ws.Range("AK10").Value = 19
$ws.Range("AN10").Value = 5.5
$ws.Range("AR10").Value = 101
$ws.Range("AU10").Value = 8.5

# Row 17
$ws.Range("G17").Value = 2.9
$ws.Range("I17").Value = 2.45
$ws.Range("M17").Value = 1.08
$ws.Range("N17").Value = 8
$ws.Range("O17").Value = 1.36
$ws.Range("P17").Value = 3
$ws.Range("Q17").Value = 2.2
$ws.Range("R17").Value = 1.65
$ws.Range("AH17").Value = 11
$ws.Range("AV17").Value = 67
$ws.Range("AW17").Value = 4.33

# Row 19
$ws.Range("M19").Value = 1.06
$ws.Range("N19").Value = 10

# Row 20
$ws.Range("G20").Value = 2.7
$ws.Range("H20").Value = 3.1
$ws.Range("I20").Value = 2.52
$ws.Range("J20").Value = 3.3
$ws.Range("K20").Value = 2
$ws.Range("M20").Value = 1.05
$ws.Range("N20").Value = 6
$ws.Range("O20").Value = 1.44
$ws.Range("P20").Value = 2.42
$ws.Range("Q20").Value = 2.25
$ws.Range("R20").Value = 1.5
$ws.Range("S20").Value = 1.45
$ws.Range("T20").Value = 2.37
$ws.Range("U20").Value = 1.98
$ws.Range("V20").Value = 1.65
$ws.Range("W20").Value = 6.9
$ws.Range("X20").Value = 12
$ws.Range("Y20").Value = 10.75
$ws.Range("Z20").Value = 30
$ws.Range("AA20").Value = 27
$ws.Range("AB20").Value = 45
$ws.Range("AC20").Value = 7.2
$ws.Range("AD20").Value = 6.1
$ws.Range("AE20").Value = 17.5
$ws.Range("AF20").Value = 110
$ws.Range("AG20").Value = 6.7
$ws.Range("AH20").Value = 11
$ws.Range("AK20").Value = 25
$ws.Range("AL20").Value = 45
$ws.Range("AN20").Value = 4.4
$ws.Range("AO20").Value = 14.5
$ws.Range("AP20").Value = 26
$ws.Range("AQ20").Value = 70
$ws.Range("AR20").Value = 120
$ws.Range("AS20").Value = 400
$ws.Range("AT20").Value = 2.32
$ws.Range("AU20").Value = 7.7
$ws.Range("AV20").Value = 80
$ws.Range("AW20").Value = 4.25
$ws.Range("AY20").Value = 25
$ws.Range("BA20").Value = 110
$ws.Range("BB20").Value = 400
